$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.223245666666665
$ws.Range("H2").Value = 24.669737
$ws.Range("I2").Value = 0.003010099901484358
$ws.Range("J2").Value = 0.003010099901484359
$ws.Range("M2").Value = 0.74396
$ws.Range("N2").Value = 2.23188
$ws.Range("O2").Value = 0.006259003216804254
$ws.Range("P2").Value = 0.006259003216804255
$ws.Range("Q2").Value = 6.117765846173332
$ws.Range("R2").Value = 55.05989261555999
$ws.Range("S2").Value = 0.00001884022496629276
$ws.Range("T2").Value = 0.00001884022496629277
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.223245666666665
$ws.Range("H3").Value = 24.669737
$ws.Range("I3").Value = 0.003010099901484358
$ws.Range("J3").Value = 0.003010099901484359
$ws.Range("M3").Value = 88.14978533333333
$ws.Range("O3").Value = 0.7416121699579786
$ws.Range("P3").Value = 0.7416121699579786
$ws.Range("Q3").Value = 724.8773402599301
$ws.Range("R3").Value = 6523.89606233937
$ws.Range("S3").Value = 0.002232326719730112
$ws.Range("T3").Value = 0.002232326719730113
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.223245666666665
$ws.Range("H4").Value = 24.669737
$ws.Range("I4").Value = 0.003010099901484358
$ws.Range("J4").Value = 0.003010099901484359
$ws.Range("M4").Value = 29.76859933333333
$ws.Range("N4").Value = 89.305798
$ws.Range("O4").Value = 0.2504459365921425
$ws.Range("P4").Value = 0.2504459365921425
$ws.Range("Q4").Value = 244.7945054705695
$ws.Range("R4").Value = 2203.150549235126
$ws.Range("S4").Value = 0.0007538672890631657
$ws.Range("T4").Value = 0.000753867289063166
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.223245666666665
$ws.Range("H5").Value = 24.669737
$ws.Range("I5").Value = 0.003010099901484358
$ws.Range("J5").Value = 0.003010099901484359
$ws.Range("M5").Value = 0.2000323333333334
$ws.Range("N5").Value = 0.6000970000000001
$ws.Range("O5").Value = 0.00168289023307462
$ws.Range("P5").Value = 0.00168289023307462
$ws.Range("Q5").Value = 1.644915018276556
$ws.Range("R5").Value = 14.804235164489
$ws.Range("S5").Value = 0.000005065667724786902
$ws.Range("T5").Value = 0.000005065667724786904
$ws.Range("I6").Value = 0.9908672065823976
$ws.Range("J6").Value = 0.9908672065823977
$ws.Range("M6").Value = 0.74396
$ws.Range("N6").Value = 2.23188
$ws.Range("O6").Value = 0.006259003216804254
$ws.Range("P6").Value = 0.006259003216804255
$ws.Range("Q6").Value = 2013.851284980173
$ws.Range("R6").Value = 18124.66156482156
$ws.Range("S6").Value = 0.006201841033425072
$ws.Range("T6").Value = 0.006201841033425073
$ws.Range("I7").Value = 0.9908672065823976
$ws.Range("J7").Value = 0.9908672065823977
$ws.Range("M7").Value = 88.14978533333333
$ws.Range("O7").Value = 0.7416121699579786
$ws.Range("P7").Value = 0.7416121699579786
$ws.Range("Q7").Value = 238615.7299643257
$ws.Range("S7").Value = 0.7348391792137725
$ws.Range("T7").Value = 0.7348391792137726
$ws.Range("I8").Value = 0.9908672065823976
$ws.Range("J8").Value = 0.9908672065823977
$ws.Range("M8").Value = 29.76859933333333
$ws.Range("N8").Value = 89.305798
$ws.Range("O8").Value = 0.2504459365921425
$ws.Range("P8").Value = 0.2504459365921425
$ws.Range("Q8").Value = 80581.66033051946
$ws.Range("R8").Value = 725234.9429746752
$ws.Range("S8").Value = 0.2481586655909685
$ws.Range("T8").Value = 0.2481586655909685
$ws.Range("I9").Value = 0.9908672065823976
$ws.Range("J9").Value = 0.9908672065823977
$ws.Range("M9").Value = 0.2000323333333334
$ws.Range("N9").Value = 0.6000970000000001
$ws.Range("O9").Value = 0.00168289023307462
$ws.Range("P9").Value = 0.00168289023307462
$ws.Range("Q9").Value = 541.4745033616267
$ws.Range("R9").Value = 4873.27053025464
$ws.Range("S9").Value = 0.001667520744231449
$ws.Range("T9").Value = 0.001667520744231449
$ws.Range("G10").Value = 14.14340733333333
$ws.Range("H10").Value = 42.430222
$ws.Range("I10").Value = 0.005177161275053701
$ws.Range("J10").Value = 0.005177161275053702
$ws.Range("M10").Value = 0.74396
$ws.Range("N10").Value = 2.23188
$ws.Range("O10").Value = 0.006259003216804254
$ws.Range("P10").Value = 0.006259003216804255
$ws.Range("Q10").Value = 10.52212931970667
$ws.Range("R10").Value = 94.69916387735999
$ws.Range("S10").Value = 0.00003240386907447553
$ws.Range("T10").Value = 0.00003240386907447554
$ws.Range("G11").Value = 14.14340733333333
$ws.Range("H11").Value = 42.430222
$ws.Range("I11").Value = 0.005177161275053701
$ws.Range("J11").Value = 0.005177161275053702
$ws.Range("M11").Value = 88.14978533333333
$ws.Range("O11").Value = 0.7416121699579786
$ws.Range("P11").Value = 0.7416121699579786
$ws.Range("Q11").Value = 1246.738320315226
$ws.Range("R11").Value = 11220.64488283703
$ws.Range("S11").Value = 0.00383944580741499
$ws.Range("T11").Value = 0.003839445807414991
$ws.Range("G12").Value = 14.14340733333333
$ws.Range("H12").Value = 42.430222
$ws.Range("I12").Value = 0.005177161275053701
$ws.Range("J12").Value = 0.005177161275053702
$ws.Range("M12").Value = 29.76859933333333
$ws.Range("N12").Value = 89.305798
$ws.Range("O12").Value = 0.2504459365921425
$ws.Range("P12").Value = 0.2504459365921425
$ws.Range("Q12").Value = 421.0294261141285
$ws.Range("R12").Value = 3789.264835027156
$ws.Range("S12").Value = 0.001296599004419395
$ws.Range("T12").Value = 0.001296599004419395
$ws.Range("G13").Value = 14.14340733333333
$ws.Range("H13").Value = 42.430222
$ws.Range("I13").Value = 0.005177161275053701
$ws.Range("J13").Value = 0.005177161275053702
$ws.Range("M13").Value = 0.2000323333333334
$ws.Range("N13").Value = 0.6000970000000001
$ws.Range("O13").Value = 0.00168289023307462
$ws.Range("P13").Value = 0.00168289023307462
$ws.Range("Q13").Value = 2.829138770170445
$ws.Range("R13").Value = 25.46224893153401
$ws.Range("S13").Value = 0.00000871259414484002
$ws.Range("T13").Value = 0.000008712594144840022
$ws.Range("G14").Value = 2.583085
$ws.Range("H14").Value = 7.749255
$ws.Range("I14").Value = 0.0009455322410643118
$ws.Range("J14").Value = 0.0009455322410643119
$ws.Range("M14").Value = 0.74396
$ws.Range("N14").Value = 2.23188
$ws.Range("O14").Value = 0.006259003216804254
$ws.Range("P14").Value = 0.006259003216804255
$ws.Range("Q14").Value = 1.9217119166
$ws.Range("R14").Value = 17.2954072494
$ws.Range("S14").Value = 0.000005918089338413663
$ws.Range("T14").Value = 0.000005918089338413664
$ws.Range("G15").Value = 2.583085
$ws.Range("H15").Value = 7.749255
$ws.Range("I15").Value = 0.0009455322410643118
$ws.Range("J15").Value = 0.0009455322410643119
$ws.Range("M15").Value = 88.14978533333333
$ws.Range("O15").Value = 0.7416121699579786
$ws.Range("P15").Value = 0.7416121699579786
$ws.Range("Q15").Value = 227.6983882477533
$ws.Range("R15").Value = 2049.285494229779
$ws.Range("S15").Value = 0.0007012182170609347
$ws.Range("T15").Value = 0.0007012182170609349
$ws.Range("G16").Value = 2.583085
$ws.Range("H16").Value = 7.749255
$ws.Range("I16").Value = 0.0009455322410643118
$ws.Range("J16").Value = 0.0009455322410643119
$ws.Range("M16").Value = 29.76859933333333
$ws.Range("N16").Value = 89.305798
$ws.Range("O16").Value = 0.2504459365921425
$ws.Range("P16").Value = 0.2504459365921425
$ws.Range("Q16").Value = 76.89482240894333
$ws.Range("R16").Value = 692.05340168049
$ws.Range("S16").Value = 0.000236804707691419
$ws.Range("T16").Value = 0.000236804707691419
$ws.Range("G17").Value = 2.583085
$ws.Range("H17").Value = 7.749255
$ws.Range("I17").Value = 0.0009455322410643118
$ws.Range("J17").Value = 0.0009455322410643119
$ws.Range("M17").Value = 0.2000323333333334
$ws.Range("N17").Value = 0.6000970000000001
$ws.Range("O17").Value = 0.00168289023307462
$ws.Range("P17").Value = 0.00168289023307462
$ws.Range("Q17").Value = 0.5167005197483334
$ws.Range("R17").Value = 4.650304677735001
$ws.Range("S17").Value = 0.000001591226973544288
$ws.Range("T17").Value = 0.000001591226973544288
